$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 5854109
$ws.Cells.Item(86, 9).Value = 7970
$ws.Cells.Item(86, 11).Value = 7970
$ws.Cells.Item(86, 13).Value = -6847
$ws.Cells.Item(87, 8).Value = 20000
$ws.Cells.Item(87, 10).Value = 20000
$ws.Cells.Item(87, 12).Value = 20000
$ws.Cells.Item(87, 14).Value = -22496
$ws.Cells.Item(89, 8).Value = 5854109
$ws.Cells.Item(89, 9).Value = 7970
$ws.Cells.Item(89, 11).Value = 39850
$ws.Cells.Item(89, 13).Value = -34234
$ws.Cells.Item(90, 8).Value = 20000
$ws.Cells.Item(90, 10).Value = 20000
$ws.Cells.Item(90, 12).Value = 60000
$ws.Cells.Item(90, 14).Value = -72480
$ws.Cells.Item(129, 8).Value = 743810.9399999999
$ws.Cells.Item(129, 9).Value = 891593.1
$ws.Cells.Item(129, 11).Value = 2674779.3
$ws.Cells.Item(129, 13).Value = -2669779.3
$ws.Cells.Item(134, 8).Value = 80000
$ws.Cells.Item(134, 10).Value = 80000
$ws.Cells.Item(134, 12).Value = 80000
$ws.Cells.Item(134, 14).Value = -90140
$ws.Cells.Item(137, 8).Value = 1646.4445
$ws.Cells.Item(137, 9).Value = 1031.7894
$ws.Cells.Item(137, 11).Value = 3095.3682
$ws.Cells.Item(137, 13).Value = -545.3681999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8628.137000000001
$ws.Cells.Item(32, 9).Value = 6296.356
$ws.Cells.Item(32, 10).Value = 18454.928
$ws.Cells.Item(32, 11).Value = 6296.356
$ws.Cells.Item(32, 12).Value = 18454.928
$ws.Cells.Item(32, 13).Value = -6009.356
$ws.Cells.Item(32, 14).Value = -19028.928
$ws.Cells.Item(97, 8).Value = 5849870.5
$ws.Cells.Item(97, 9).Value = 1829.6364
$ws.Cells.Item(97, 11).Value = 1829.6364
$ws.Cells.Item(97, 13).Value = -1333.6364
$ws.Cells.Item(102, 8).Value = 17546112
$ws.Cells.Item(102, 9).Value = 2143.125
$ws.Cells.Item(102, 11).Value = 2143.125
$ws.Cells.Item(102, 13).Value = -521.125
$ws.Cells.Item(132, 8).Value = 4429.5415
$ws.Cells.Item(132, 9).Value = 3832.875
$ws.Cells.Item(132, 11).Value = 11498.625
$ws.Cells.Item(132, 13).Value = -8968.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 33749.75
$ws.Cells.Item(76, 9).Value = 38499.5
$ws.Cells.Item(76, 10).Value = 29000
$ws.Cells.Item(76, 11).Value = 38499.5
$ws.Cells.Item(76, 12).Value = 29000
$ws.Cells.Item(76, 13).Value = -38184.5
$ws.Cells.Item(76, 14).Value = -29630
$ws.Cells.Item(79, 8).Value = 33749.75
$ws.Cells.Item(79, 9).Value = 38499.5
$ws.Cells.Item(79, 10).Value = 29000
$ws.Cells.Item(79, 11).Value = 38499.5
$ws.Cells.Item(79, 12).Value = 29000
$ws.Cells.Item(79, 13).Value = -37407.5
$ws.Cells.Item(79, 14).Value = -31184
$ws.Cells.Item(134, 8).Value = 5141.0347
$ws.Cells.Item(134, 9).Value = 4454.3887
$ws.Cells.Item(134, 10).Value = 6264.636
$ws.Cells.Item(134, 11).Value = 13363.1661
$ws.Cells.Item(134, 12).Value = 18793.908
$ws.Cells.Item(134, 13).Value = -10828.1661
$ws.Cells.Item(134, 14).Value = -23863.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 336.3
$ws.Cells.Item(7, 9).Value = 30
$ws.Cells.Item(7, 10).Value = 438.4
$ws.Cells.Item(7, 11).Value = 30
$ws.Cells.Item(7, 12).Value = 438.4
$ws.Cells.Item(7, 13).Value = 83
$ws.Cells.Item(7, 14).Value = -664.4
$ws.Cells.Item(31, 9).Value = 1505
$ws.Cells.Item(31, 10).Value = 4719
$ws.Cells.Item(31, 11).Value = 1505
$ws.Cells.Item(31, 12).Value = 4719
$ws.Cells.Item(31, 13).Value = -1210
$ws.Cells.Item(31, 14).Value = -5309
$ws.Cells.Item(34, 9).Value = 1505
$ws.Cells.Item(34, 10).Value = 4719
$ws.Cells.Item(34, 11).Value = 1505
$ws.Cells.Item(34, 12).Value = 4719
$ws.Cells.Item(34, 13).Value = -1303
$ws.Cells.Item(34, 14).Value = -5123
$ws.Cells.Item(62, 8).Value = 9096.75
$ws.Cells.Item(62, 9).Value = 10346.75
$ws.Cells.Item(62, 10).Value = 7846.75
$ws.Cells.Item(62, 11).Value = 10346.75
$ws.Cells.Item(62, 12).Value = 7846.75
$ws.Cells.Item(62, 13).Value = -9722.75
$ws.Cells.Item(62, 14).Value = -9094.75
$ws.Cells.Item(65, 8).Value = 9096.75
$ws.Cells.Item(65, 9).Value = 10346.75
$ws.Cells.Item(65, 10).Value = 7846.75
$ws.Cells.Item(65, 11).Value = 51733.75
$ws.Cells.Item(65, 12).Value = 39233.75
$ws.Cells.Item(65, 13).Value = -48613.75
$ws.Cells.Item(65, 14).Value = -45473.75
$ws.Cells.Item(132, 8).Value = 2158.8076
$ws.Cells.Item(132, 9).Value = 1907.6
$ws.Cells.Item(132, 11).Value = 5722.799999999999
$ws.Cells.Item(132, 13).Value = -3192.799999999999
$ws.Cells.Item(134, 8).Value = 4612.9473
$ws.Cells.Item(134, 9).Value = 3537.3333
$ws.Cells.Item(134, 11).Value = 10611.9999
$ws.Cells.Item(134, 13).Value = -8076.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 32820088
$ws.Cells.Item(4, 9).Value = 37773948
$ws.Cells.Item(4, 11).Value = 113321844
$ws.Cells.Item(4, 13).Value = -113321732
$ws.Cells.Item(55, 8).Value = 913155.25
$ws.Cells.Item(55, 9).Value = 3900
$ws.Cells.Item(55, 10).Value = 1004080.8
$ws.Cells.Item(55, 11).Value = 11700
$ws.Cells.Item(55, 12).Value = 3012242.4
$ws.Cells.Item(55, 13).Value = -11523
$ws.Cells.Item(55, 14).Value = -3012596.4
$ws.Cells.Item(99, 8).Value = 61289.6
$ws.Cells.Item(99, 10).Value = 74862
$ws.Cells.Item(99, 12).Value = 224586
$ws.Cells.Item(99, 14).Value = -229078
$ws.Cells.Item(115, 8).Value = 2549
$ws.Cells.Item(115, 9).Value = 971.5
$ws.Cells.Item(115, 10).Value = 3180
$ws.Cells.Item(115, 11).Value = 2914.5
$ws.Cells.Item(115, 12).Value = 9540
$ws.Cells.Item(115, 13).Value = -1739.5
$ws.Cells.Item(115, 14).Value = -11890
$ws.Cells.Item(121, 8).Value = 1870.1818
$ws.Cells.Item(121, 9).Value = 239.75
$ws.Cells.Item(121, 10).Value = 2391.92
$ws.Cells.Item(121, 11).Value = 719.25
$ws.Cells.Item(121, 12).Value = 7175.76
$ws.Cells.Item(121, 13).Value = 590.75
$ws.Cells.Item(121, 14).Value = -9795.76

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 318.65
$ws.Cells.Item(2, 10).Value = 349.44446
$ws.Cells.Item(2, 12).Value = 349.44446
$ws.Cells.Item(2, 14).Value = -575.4444599999999
$ws.Cells.Item(70, 8).Value = 191454.33
$ws.Cells.Item(70, 9).Value = 225745.2
$ws.Cells.Item(70, 11).Value = 225745.2
$ws.Cells.Item(70, 13).Value = -225475.2
$ws.Cells.Item(73, 8).Value = 191454.33
$ws.Cells.Item(73, 9).Value = 225745.2
$ws.Cells.Item(73, 11).Value = 225745.2
$ws.Cells.Item(73, 13).Value = -224809.2
$ws.Cells.Item(80, 8).Value = 27872878
$ws.Cells.Item(80, 9).Value = 187366.83
$ws.Cells.Item(80, 10).Value = 55558390
$ws.Cells.Item(80, 11).Value = 187366.83
$ws.Cells.Item(80, 12).Value = 55558390
$ws.Cells.Item(80, 13).Value = -186368.83
$ws.Cells.Item(80, 14).Value = -55560386
$ws.Cells.Item(83, 8).Value = 27872878
$ws.Cells.Item(83, 9).Value = 187366.83
$ws.Cells.Item(83, 10).Value = 55558390
$ws.Cells.Item(83, 11).Value = 936834.1499999999
$ws.Cells.Item(83, 12).Value = 277791950
$ws.Cells.Item(83, 13).Value = -931842.1499999999
$ws.Cells.Item(83, 14).Value = -277801934
$ws.Cells.Item(132, 8).Value = 2494.2292
$ws.Cells.Item(132, 9).Value = 1990.6
$ws.Cells.Item(132, 10).Value = 3850.1538
$ws.Cells.Item(132, 11).Value = 5971.799999999999
$ws.Cells.Item(132, 12).Value = 11550.4614
$ws.Cells.Item(132, 13).Value = -3441.799999999999
$ws.Cells.Item(132, 14).Value = -16610.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 7671.875
$ws.Cells.Item(61, 10).Value = 7978
$ws.Cells.Item(61, 12).Value = 7978
$ws.Cells.Item(61, 14).Value = -8382
$ws.Cells.Item(113, 8).Value = 7671.875
$ws.Cells.Item(113, 10).Value = 7978
$ws.Cells.Item(113, 12).Value = 7978
$ws.Cells.Item(113, 14).Value = -12318
$ws.Cells.Item(132, 8).Value = 4353.92
$ws.Cells.Item(132, 9).Value = 3492.875
$ws.Cells.Item(132, 11).Value = 10478.625
$ws.Cells.Item(132, 13).Value = -7948.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(117, 8).Value = 180000
$ws.Cells.Item(117, 10).Value = 180000
$ws.Cells.Item(117, 12).Value = 180000
$ws.Cells.Item(117, 14).Value = -189178
$ws.Cells.Item(124, 8).Value = 399966.66
$ws.Cells.Item(124, 10).Value = 399966.66
$ws.Cells.Item(124, 12).Value = 399966.66
$ws.Cells.Item(124, 14).Value = -409786.66
$ws.Cells.Item(132, 8).Value = 1827.242
$ws.Cells.Item(132, 9).Value = 1391.6052
$ws.Cells.Item(132, 10).Value = 2517
$ws.Cells.Item(132, 11).Value = 4174.8156
$ws.Cells.Item(132, 12).Value = 7551
$ws.Cells.Item(132, 13).Value = -1644.8156
$ws.Cells.Item(132, 14).Value = -12611
$ws.Cells.Item(136, 8).Value = 5206.1665
$ws.Cells.Item(136, 9).Value = 1749.5
$ws.Cells.Item(136, 10).Value = 5897.5
$ws.Cells.Item(136, 11).Value = 5248.5
$ws.Cells.Item(136, 12).Value = 17692.5
$ws.Cells.Item(136, 13).Value = -2698.5
$ws.Cells.Item(136, 14).Value = -22792.5
